# Add daily power records
# - Fix C6 / C10 "End Time" typos (were AM, should be the near-midnight PM value)
# - Extend the comforter_cda_table (and the underlying data block) with a new
#   day (row 57, date 43381) following the same pattern as the other "no
#   reading" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Correct two mis-keyed End Time values -------------------------------
$ws.Range("C6").Value = 0.99930555555555556
$ws.Range("C10").Value = 0.99930555555555556

# --- Rows 52/53/56 were missing explicit Start/End Time zeros; Excel
#     fills these in as part of extending the calculated-column table. ----
$ws.Range("B52").Value = 0
$ws.Range("C52").Value = 0
$ws.Range("B53").Value = 0
$ws.Range("C53").Value = 0
$ws.Range("B56").Value = 0
$ws.Range("C56").Value = 0

# --- Re-fill the calculated columns D51:F56 so Excel collapses them back
#     into shared formulas (matching how Excel represents a filled-down
#     formula block). ------------------------------------------------------
$ws.Range("D51:D56").Formula = "=(C51-B51)* 1440"
$ws.Range("E51:E56").Formula = "=IF(C51>B51, (C51-B51)*1440, (B51-C51)*1440)"
$ws.Range("F51:F56").Formula = "=ABS((C51-B51)*1440)"

# --- Add the new daily record (row 57) ------------------------------------
$ws.Range("A57").Value = 43381
$ws.Range("D57").Formula = "=(C57-B57)* 1440"
$ws.Range("E57").Formula = "=IF(C57>B57, (C57-B57)*1440, (B57-C57)*1440)"
$ws.Range("F57").Formula = "=ABS((C57-B57)*1440)"

# --- Keep the table definition (ref + autofilter) in sync with the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F57"))

# --- Mirror the view state move (Excel scrolls/selects the new last row)
$ws.Range("B57").Select()
$ws.Application.ActiveWindow.ScrollRow = 46
